$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2772.5
$ws.Range("J40").Value = 2827.4285
$ws.Range("L40").Value = 2827.4285
$ws.Range("N40").Value = -3177.4285

$ws.Range("H62").Value = 125001250
$ws.Range("I62").Value = 250000000
$ws.Range("K62").Value = 250000000
$ws.Range("M62").Value = -249999376

$ws.Range("H65").Value = 125001250
$ws.Range("I65").Value = 250000000
$ws.Range("K65").Value = 1250000000
$ws.Range("M65").Value = -1249996880

$ws.Range("H76").Value = 166671660
$ws.Range("I76").Value = 333336130
$ws.Range("J76").Value = 7185
$ws.Range("K76").Value = 333336130
$ws.Range("L76").Value = 7185
$ws.Range("M76").Value = -333335815
$ws.Range("N76").Value = -7815

$ws.Range("H79").Value = 166671660
$ws.Range("I79").Value = 333336130
$ws.Range("J79").Value = 7185
$ws.Range("K79").Value = 333336130
$ws.Range("L79").Value = 7185
$ws.Range("M79").Value = -333335038
$ws.Range("N79").Value = -9369

$ws.Range("H107").Value = 2102.8572
$ws.Range("I107").Value = 2225.4583
$ws.Range("K107").Value = 2225.4583
$ws.Range("M107").Value = -305.4582999999998

$ws.Range("H132").Value = 3949.1614
$ws.Range("I132").Value = 3793.4814
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 11380.4442
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -8850.4442
$ws.Range("N132").Value = -20060

$ws.Range("H135").Value = 28706.182
$ws.Range("I135").Value = 1953.4
$ws.Range("J135").Value = 51000.168
$ws.Range("K135").Value = 17580.6
$ws.Range("L135").Value = 459001.512
$ws.Range("M135").Value = -15045.6
$ws.Range("N135").Value = -464071.512

$ws.Range("H138").Value = 2911.5593
$ws.Range("I138").Value = 2310.2068
$ws.Range("K138").Value = 6930.6204
$ws.Range("M138").Value = -1790.6204

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1884.5555
$ws.Range("I2").Value = 1745.125
$ws.Range("K2").Value = 1745.125
$ws.Range("M2").Value = -1632.125

$ws.Range("H5").Value = 55
$ws.Range("I5").Value = 46
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 46
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 66
$ws.Range("N5").Value = -324

$ws.Range("H61").Value = 5848.9575
$ws.Range("I61").Value = 4544.5356
$ws.Range("J61").Value = 7771.263
$ws.Range("K61").Value = 4544.5356
$ws.Range("L61").Value = 7771.263
$ws.Range("M61").Value = -4332.5356
$ws.Range("N61").Value = -8195.262999999999

$ws.Range("H63").Value = 2946
$ws.Range("I63").Value = 2946
$ws.Range("K63").Value = 2946
$ws.Range("M63").Value = -2260

$ws.Range("H66").Value = 2946
$ws.Range("I66").Value = 2946
$ws.Range("K66").Value = 14730
$ws.Range("M66").Value = -11298

$ws.Range("H74").Value = 2104.425
$ws.Range("I74").Value = 1578.8823
$ws.Range("J74").Value = 5082.5
$ws.Range("K74").Value = 1578.8823
$ws.Range("L74").Value = 5082.5
$ws.Range("M74").Value = -704.8823
$ws.Range("N74").Value = -6830.5

$ws.Range("H77").Value = 2104.425
$ws.Range("I77").Value = 1578.8823
$ws.Range("J77").Value = 5082.5
$ws.Range("K77").Value = 7894.4115
$ws.Range("L77").Value = 25412.5
$ws.Range("M77").Value = -3526.4115
$ws.Range("N77").Value = -34148.5

$ws.Range("H116").Value = 1884.5555
$ws.Range("I116").Value = 1745.125
$ws.Range("K116").Value = 1745.125
$ws.Range("M116").Value = 548.875

$ws.Range("H122").Value = 2339.111
$ws.Range("I122").Value = 1360
$ws.Range("K122").Value = 4080
$ws.Range("M122").Value = -1630

$ws.Range("H136").Value = 5848.9575
$ws.Range("I136").Value = 4544.5356
$ws.Range("J136").Value = 7771.263
$ws.Range("K136").Value = 13633.6068
$ws.Range("L136").Value = 23313.789
$ws.Range("M136").Value = -11083.6068
$ws.Range("N136").Value = -28413.789

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1884.5555
$ws.Range("I3").Value = 1745.125
$ws.Range("K3").Value = 1745.125
$ws.Range("M3").Value = -1631.125

$ws.Range("H4").Value = 55
$ws.Range("I4").Value = 46
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 46
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 69
$ws.Range("N4").Value = -330

$ws.Range("H22").Value = 1250491.8
$ws.Range("I22").Value = 276.2143
$ws.Range("K22").Value = 276.2143
$ws.Range("M22").Value = -103.2143

$ws.Range("H134").Value = 5484.9023
$ws.Range("I134").Value = 4919.125
$ws.Range("J134").Value = 7496.5557
$ws.Range("K134").Value = 14757.375
$ws.Range("L134").Value = 22489.6671
$ws.Range("M134").Value = -12222.375
$ws.Range("N134").Value = -27559.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1266.6666
$ws.Range("I10").Value = 1700
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 1700
$ws.Range("L10").Value = 400
$ws.Range("M10").Value = -1561
$ws.Range("N10").Value = -678

$ws.Range("H13").Value = 1749.5
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 2499
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 2499
$ws.Range("M13").Value = -861
$ws.Range("N13").Value = -2777

$ws.Range("H31").Value = 3217.8
$ws.Range("I31").Value = 1805.4546
$ws.Range("J31").Value = 4944
$ws.Range("K31").Value = 1805.4546
$ws.Range("L31").Value = 4944
$ws.Range("M31").Value = -1510.4546
$ws.Range("N31").Value = -5534

$ws.Range("H34").Value = 3217.8
$ws.Range("I34").Value = 1805.4546
$ws.Range("J34").Value = 4944
$ws.Range("K34").Value = 1805.4546
$ws.Range("L34").Value = 4944
$ws.Range("M34").Value = -1603.4546
$ws.Range("N34").Value = -5348

$ws.Range("H107").Value = 3304.7083
$ws.Range("I107").Value = 4008.7058
$ws.Range("K107").Value = 4008.7058
$ws.Range("M107").Value = -2088.7058

$ws.Range("H122").Value = 4586.25
$ws.Range("I122").Value = 4458.923
$ws.Range("J122").Value = 5138
$ws.Range("K122").Value = 13376.769
$ws.Range("L122").Value = 15414
$ws.Range("M122").Value = -10926.769
$ws.Range("N122").Value = -20314

$ws.Range("H132").Value = 2429.842
$ws.Range("I132").Value = 2231.5
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 6694.5
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -4164.5
$ws.Range("N132").Value = -23060

$ws.Range("H134").Value = 4984.049
$ws.Range("I134").Value = 4421.4707
$ws.Range("K134").Value = 13264.4121
$ws.Range("M134").Value = -10729.4121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1490
$ws.Range("J107").Value = 1739.1666
$ws.Range("L107").Value = 5217.4998
$ws.Range("N107").Value = -9057.4998

$ws.Range("H112").Value = 7805.6
$ws.Range("I112").Value = 5666
$ws.Range("J112").Value = 11015
$ws.Range("K112").Value = 16998
$ws.Range("L112").Value = 33045
$ws.Range("M112").Value = -15890
$ws.Range("N112").Value = -35261

$ws.Range("H113").Value = 993.75
$ws.Range("I113").Value = 753.3333
$ws.Range("J113").Value = 1715
$ws.Range("K113").Value = 2259.9999
$ws.Range("L113").Value = 5145
$ws.Range("M113").Value = -89.9998999999998
$ws.Range("N113").Value = -9485

$ws.Range("H131").Value = 23394186
$ws.Range("I131").Value = 7938562.5
$ws.Range("J131").Value = 66669930
$ws.Range("K131").Value = 23815687.5
$ws.Range("L131").Value = 200009790
$ws.Range("M131").Value = -23810647.5
$ws.Range("N131").Value = -200019870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2618.3044
$ws.Range("I80").Value = 2325.2
$ws.Range("J80").Value = 2843.7693
$ws.Range("K80").Value = 2325.2
$ws.Range("L80").Value = 2843.7693
$ws.Range("M80").Value = -1327.2
$ws.Range("N80").Value = -4839.7693

$ws.Range("H83").Value = 2618.3044
$ws.Range("I83").Value = 2325.2
$ws.Range("J83").Value = 2843.7693
$ws.Range("K83").Value = 11626
$ws.Range("L83").Value = 14218.8465
$ws.Range("M83").Value = -6634
$ws.Range("N83").Value = -24202.8465

$ws.Range("H126").Value = 4980.6924
$ws.Range("J126").Value = 6083.3335
$ws.Range("L126").Value = 18250.0005
$ws.Range("N126").Value = -23190.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2016.2759
$ws.Range("I40").Value = 1740.6666
$ws.Range("J40").Value = 3339.2
$ws.Range("K40").Value = 1740.6666
$ws.Range("L40").Value = 3339.2
$ws.Range("M40").Value = -1604.6666
$ws.Range("N40").Value = -3611.2

$ws.Range("H82").Value = 2687.6155
$ws.Range("I82").Value = 2232.7778
$ws.Range("J82").Value = 3711
$ws.Range("K82").Value = 2232.7778
$ws.Range("L82").Value = 3711
$ws.Range("M82").Value = -1871.7778
$ws.Range("N82").Value = -4433

$ws.Range("H85").Value = 2687.6155
$ws.Range("I85").Value = 2232.7778
$ws.Range("J85").Value = 3711
$ws.Range("K85").Value = 2232.7778
$ws.Range("L85").Value = 3711
$ws.Range("M85").Value = -984.7777999999998
$ws.Range("N85").Value = -6207

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8974.75

$ws.Range("H65").Value = 8974.75

$ws.Range("H122").Value = 326550.56
$ws.Range("I122").Value = 503525.34
$ws.Range("J122").Value = 4778.1816
$ws.Range("K122").Value = 1510576.02
$ws.Range("L122").Value = 14334.5448
$ws.Range("M122").Value = -1508126.02
$ws.Range("N122").Value = -19234.5448
